$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.659.02'
$ws.Range("E2").Value = '  +1.39%  '

# Row 3
$ws.Range("D3").Value = '3.727.15'
$ws.Range("E3").Value = '  -0.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.64%  '

# Row 7
$ws.Range("D7").Value = '3.723.77'
$ws.Range("E7").Value = '  -0.94%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("E9").Value = '  +1.25%  '

# Row 10
$ws.Range("E10").Value = '  +2.86%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.459'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.19%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.51%  '

# Row 14
$ws.Range("E14").Value = '  -0.49%  '

# Row 15
$ws.Range("D15").Value = '4.349.07'
$ws.Range("E15").Value = '  -1.16%  '

# Row 16
$ws.Range("D16").Value = '3.727.61'
$ws.Range("E16").Value = '  -1.05%  '

# Row 17
$ws.Range("D17").Value = '68.624.86'
$ws.Range("E17").Value = '  +1.39%  '

# Row 18
$ws.Range("E18").Value = '  +0.55%  '

# Row 19
$ws.Range("E19").Value = '  +0.28%  '

# Row 20
$ws.Range("E20").Value = '  +4.48%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '495.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.12%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +11.58%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.722'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.62%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.47%  '

# Row 25
$ws.Range("E25").Value = '  -4.37%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.17%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.72%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.08%  '

# Row 29
$ws.Range("E29").Value = '  +0.00%  '

# Row 30
$ws.Range("E30").Value = '  -0.50%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.40%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.54%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.70%  '

# Row 34
$ws.Range("D34").Value = '3.871.78'
$ws.Range("E34").Value = '  -0.85%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.108'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.92%  '

# Row 36
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '3.658.39'
$ws.Range("E36").Value = '  -1.10%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.19%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.36%  '

# Row 40
$ws.Range("E40").Value = '  -1.68%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.325'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.85%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '433.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.07%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.04%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.16%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.91%  '

# Row 47
$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.17%  '

# Row 48
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.01%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.88%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0353'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.92%  '

# Row 51
$ws.Range("D51").Value = '2.742.90'
$ws.Range("E51").Value = '  -3.03%  '
